# chore: update config and sync Excel file to repository
#
# Ferrum Decor Studio product-tracking sheet:
#  - The "PURE BRASS Personalized Letter Box" price ($949.00 USD) listing
#    had its "From " prefix removed (it is no longer a starting-at price).
#    This shared-string text is reused by three rows (A65, A74, A158); the
#    last of those (A158) also still carries the trailing pencil emoji.
#  - The "Pure brass personalized mailbox" listing's price ($395.00 USD)
#    also had its "From " prefix removed, dropping the trailing emoji too.
#  - The active selection moved from A144 to A64.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Price text corrections -------------------------------------------------
$ws.Range("A65").Value  = "$949.00 USD"
$ws.Range("A74").Value  = "$949.00 USD"
$ws.Range("A149").Value = "$395.00 USD"
$ws.Range("A158").Value = "$949.00 USD✏️"

# --- Selection / view update -------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A64").Select() | Out-Null
